# Nina Williams Tekken 8 frame-data sheet: the "Block" frame-advantage
# column (E) had its sign convention reversed (e.g. -1 -> 1, 3 -> -3).
# Flip the sign of every numeric value in column E for the data rows,
# leaving blank cells (headers/notes rows with no block value) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)   # Column E = "Block"
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $cell.Value2 = -1 * $val
    }
}
